$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Title heading and matching bold footer text (identical string, both occurrences replaced)
Replace-Text "Play 'Ancient Egypt' Free - Exciting Egyptian Slot Machine" "Play Ancient Egypt Free - Exciting Ancient Egyptian Themed Slot Game"

# "What we like" bullet list
Replace-Text "Exciting Egyptian theme" "Prevalent and popular theme of ancient Egypt"
Replace-Text "Solid potential payouts" "Beautiful graphics and well-crafted design"
Replace-Text "Well-crafted graphics" "Solid potential payouts and exciting special features"

# "What we don't like" bullet list
Replace-Text "Limited special features" "Limited number of paylines"
Replace-Text "Traditional playing cards symbols" "Limited number of special features"

# Closing italic summary text
Replace-Text "Explore solid payout potential in this Egyptian themed game with well-crafted graphics. Play 'Ancient Egypt' for free and check it out today." "Read our review of Ancient Egypt, an enjoyable and potentially rewarding slot game. Play for free!"
